$p = $ppt.ActivePresentation

# Slide 1: Title text is "Header" + " " + "with" + " " + "inline code" (Courier).
# Consolidate the first four plain (unformatted) runs -> "Header with ",
# leaving the differently-formatted "inline code" run as-is.
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(1).TextFrame.TextRange
$c1 = $tr1.Characters(1, 12)
$c1.Text = "............"
$c1 = $tr1.Characters(1, 12)
$c1.Text = "Header with "

# Slide 2: Title text is "Syntax" + " " + "highlighting" (all plain runs).
# Consolidate into a single run "Syntax highlighting".
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(1).TextFrame.TextRange
$tr2.Text = "...................."
$tr2.Text = "Syntax highlighting"

# Slide 3: Title text is "Two" + " " + "column" + " " + "slide" (all plain runs).
# Consolidate into a single run "Two column slide".
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(1).TextFrame.TextRange
$tr3.Text = "................."
$tr3.Text = "Two column slide"
